# Daily auto-push edit: insert one new data row (2026/02/28, 土, 19, 39)
# immediately before the existing row 878, pushing every row from the old
# 878 through 919 down by one (new 879..920). This matches the upstream
# diff: dimension grows from A1:D919 to A1:D920, and the "header" rows for
# each date group (878, 882, 888, 894, 899, 906, 913, 918) shift down by
# one row while keeping their own text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at position 878; Excel shifts rows 878:919 -> 879:920
# and auto-updates the sheet's dimension (A1:D919 -> A1:D920).
$ws.Rows.Item(878).Insert()

# Column A holds the date as plain text (e.g. "2026/02/28"), not a real
# Excel date serial. Force text formatting first so the slash-separated
# value isn't auto-parsed into a date, then drop back to the default
# "Normal" style so the cell doesn't carry a stray explicit style (keeping
# it identical in shape to every other untouched data cell in the sheet).
$ws.Cells.Item(878,1).NumberFormat = "@"
$ws.Cells.Item(878,1).Value() = "2026/02/28"
$ws.Cells.Item(878,1).Style = "Normal"

# Weekday (text) and the two numeric columns for the new row.
$ws.Cells.Item(878,2).Value() = "土"
$ws.Cells.Item(878,3).Value() = 19
$ws.Cells.Item(878,4).Value() = 39
